# Insert a new data row at row 22 (pushing the existing rows 22-62 down to
# 23-63) and populate it with the new weekly price observation for
# Ají / Americana (o), Región del Maule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 22..62 down by one (Excel re-numbers their formulas /
# references automatically, same as pressing "Insert Sheet Rows" above row 22).
$ws.Rows.Item(22).Insert()

# Populate the newly-inserted (blank) row 22 with the new record.
$ws.Cells.Item(22, 1).Value2  = 7
$ws.Cells.Item(22, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(22, 3).Value2  = "Ñuble"
$ws.Cells.Item(22, 4).Value2  = 44592
$ws.Cells.Item(22, 5).Value2  = 16
$ws.Cells.Item(22, 6).Value2  = 100112021
$ws.Cells.Item(22, 7).Value2  = "Ají"
$ws.Cells.Item(22, 8).Value2  = "Americana (o)"
$ws.Cells.Item(22, 9).Value2  = "Primera"
$ws.Cells.Item(22, 10).Value2 = 60
$ws.Cells.Item(22, 11).Value2 = 10000
$ws.Cells.Item(22, 12).Value2 = 10500
$ws.Cells.Item(22, 13).Value2 = 10250
$ws.Cells.Item(22, 14).Value2 = "`$/caja 15 kilos"
$ws.Cells.Item(22, 15).Value2 = "Región del Maule"
$ws.Cells.Item(22, 16).Value2 = 683
$ws.Cells.Item(22, 17).Value2 = 15
$ws.Cells.Item(22, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same date/time number format as the rest
# of column D (style index "2" => numFmtId 165, YYYY-MM-DD HH:MM:SS).
$ws.Cells.Item(22, 4).NumberFormat = $ws.Cells.Item(23, 4).NumberFormat
